$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1695.1666
$ws.Range("J18").Value = 3489.5
$ws.Range("L18").Value = 3489.5
$ws.Range("N18").Value = -4057.5

$ws.Range("H20").Value = 11933
$ws.Range("I20").Value = 5399.5
$ws.Range("J20").Value = 25000
$ws.Range("K20").Value = 5399.5
$ws.Range("L20").Value = 25000
$ws.Range("M20").Value = -5169.5
$ws.Range("N20").Value = -25460

$ws.Range("H35").Value = 11933
$ws.Range("I35").Value = 5399.5
$ws.Range("J35").Value = 25000
$ws.Range("K35").Value = 5399.5
$ws.Range("L35").Value = 25000
$ws.Range("M35").Value = -5020.5
$ws.Range("N35").Value = -25758

$ws.Range("H41").Value = 996.3333
$ws.Range("I41").Value = 1335.625
$ws.Range("K41").Value = 1335.625
$ws.Range("M41").Value = -895.625

$ws.Range("H64").Value = 29872.25
$ws.Range("I64").Value = 6745
$ws.Range("K64").Value = 6745
$ws.Range("M64").Value = -6497

$ws.Range("H67").Value = 29872.25
$ws.Range("I67").Value = 6745
$ws.Range("K67").Value = 6745
$ws.Range("M67").Value = -5887

$ws.Range("H92").Value = 2992.1428
$ws.Range("I92").Value = 3109.2
$ws.Range("K92").Value = 3109.2
$ws.Range("M92").Value = -1861.2

$ws.Range("H94").Value = 637.5
$ws.Range("I94").Value = 637.5
$ws.Range("K94").Value = 637.5
$ws.Range("M94").Value = -186.5

$ws.Range("H106").Value = 5707.4165
$ws.Range("I106").Value = 3697.7144
$ws.Range("K106").Value = 3697.7144
$ws.Range("M106").Value = -3066.7144

$ws.Range("H112").Value = 1319731.4
$ws.Range("J112").Value = 2506349.8
$ws.Range("L112").Value = 7519049.399999999
$ws.Range("N112").Value = -7521265.399999999

$ws.Range("H129").Value = 4596.143
$ws.Range("I129").Value = 1141.7142
$ws.Range("J129").Value = 6323.357
$ws.Range("K129").Value = 3425.1426
$ws.Range("L129").Value = 18970.071
$ws.Range("M129").Value = 1574.8574
$ws.Range("N129").Value = -28970.071

$ws.Range("H135").Value = 4236.68
$ws.Range("J135").Value = 9163.833000000001
$ws.Range("L135").Value = 82474.497
$ws.Range("N135").Value = -87544.497

$ws.Range("H138").Value = 5923.881
$ws.Range("I138").Value = 3203.0588
$ws.Range("J138").Value = 7774.04
$ws.Range("K138").Value = 9609.1764
$ws.Range("L138").Value = 23322.12
$ws.Range("M138").Value = -4469.1764
$ws.Range("N138").Value = -33602.12

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5522.625
$ws.Range("I2").Value = 5670.143
$ws.Range("K2").Value = 5670.143
$ws.Range("M2").Value = -5557.143

$ws.Range("H32").Value = 3570.1538
$ws.Range("I32").Value = 3145.9556
$ws.Range("J32").Value = 6297.143
$ws.Range("K32").Value = 3145.9556
$ws.Range("L32").Value = 6297.143
$ws.Range("M32").Value = -2858.9556
$ws.Range("N32").Value = -6871.143

$ws.Range("H116").Value = 5522.625
$ws.Range("I116").Value = 5670.143
$ws.Range("K116").Value = 5670.143
$ws.Range("M116").Value = -3376.143

$ws.Range("H122").Value = 2932.5757
$ws.Range("I122").Value = 2994.7188
$ws.Range("K122").Value = 8984.1564
$ws.Range("M122").Value = -6534.1564

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5026
$ws.Range("I3").Value = 5115.3335
$ws.Range("K3").Value = 5115.3335
$ws.Range("M3").Value = -5001.3335

$ws.Range("H86").Value = 2427.611
$ws.Range("I86").Value = 1602.6666
$ws.Range("J86").Value = 4077.5
$ws.Range("K86").Value = 1602.6666
$ws.Range("L86").Value = 4077.5
$ws.Range("M86").Value = -479.6666
$ws.Range("N86").Value = -6323.5

$ws.Range("H89").Value = 2427.611
$ws.Range("I89").Value = 1602.6666
$ws.Range("J89").Value = 4077.5
$ws.Range("K89").Value = 8013.333000000001
$ws.Range("L89").Value = 20387.5
$ws.Range("M89").Value = -2397.333000000001
$ws.Range("N89").Value = -31619.5

$ws.Range("H134").Value = 2529.4546
$ws.Range("I134").Value = 2425.1333
$ws.Range("K134").Value = 7275.3999
$ws.Range("M134").Value = -4740.3999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 799.3333
$ws.Range("I22").Value = 544.2
$ws.Range("K22").Value = 544.2
$ws.Range("M22").Value = -194.2

$ws.Range("H31").Value = 2560.5908
$ws.Range("I31").Value = 2621.1875
$ws.Range("K31").Value = 2621.1875
$ws.Range("M31").Value = -2326.1875

$ws.Range("H34").Value = 2560.5908
$ws.Range("I34").Value = 2621.1875
$ws.Range("K34").Value = 2621.1875
$ws.Range("M34").Value = -2419.1875

$ws.Range("H63").Value = 79268.664
$ws.Range("J63").Value = 79268.664
$ws.Range("L63").Value = 79268.664
$ws.Range("N63").Value = -80640.664

$ws.Range("H66").Value = 79268.664
$ws.Range("J66").Value = 79268.664
$ws.Range("L66").Value = 237805.992
$ws.Range("N66").Value = -244669.992

$ws.Range("H69").Value = 66749.25
$ws.Range("I69").Value = 55666
$ws.Range("J69").Value = 99999
$ws.Range("K69").Value = 55666
$ws.Range("L69").Value = 99999
$ws.Range("M69").Value = -54917
$ws.Range("N69").Value = -101497

$ws.Range("H72").Value = 66749.25
$ws.Range("I72").Value = 55666
$ws.Range("J72").Value = 99999
$ws.Range("K72").Value = 166998
$ws.Range("L72").Value = 299997
$ws.Range("M72").Value = -163254
$ws.Range("N72").Value = -307485

$ws.Range("H80").Value = 60998
$ws.Range("J80").Value = 60998
$ws.Range("L80").Value = 60998
$ws.Range("N80").Value = -63244

$ws.Range("H83").Value = 60998
$ws.Range("J83").Value = 60998
$ws.Range("L83").Value = 182994
$ws.Range("N83").Value = -194226

$ws.Range("H107").Value = 2009.619
$ws.Range("I107").Value = 1186.4286
$ws.Range("J107").Value = 2421.2144
$ws.Range("K107").Value = 1186.4286
$ws.Range("L107").Value = 2421.2144
$ws.Range("M107").Value = 733.5714
$ws.Range("N107").Value = -6261.2144

$ws.Range("H134").Value = 3601.2903
$ws.Range("I134").Value = 3665.6667
$ws.Range("K134").Value = 10997.0001
$ws.Range("M134").Value = -8462.000100000001

$ws.Range("H135").Value = 99997
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7499.4287
$ws.Range("I3").Value = 5512.231
$ws.Range("K3").Value = 16536.693
$ws.Range("M3").Value = -16424.693

$ws.Range("H121").Value = 4766189
$ws.Range("I121").Value = 684
$ws.Range("J121").Value = 6255409
$ws.Range("K121").Value = 2052
$ws.Range("L121").Value = 18766227
$ws.Range("M121").Value = -742
$ws.Range("N121").Value = -18768847

$ws.Range("H122").Value = 46142.5
$ws.Range("I122").Value = 110480
$ws.Range("K122").Value = 994320
$ws.Range("M122").Value = -991870

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 57688.332
$ws.Range("J62").Value = 57688.332
$ws.Range("L62").Value = 57688.332
$ws.Range("N62").Value = -59060.332

$ws.Range("H65").Value = 57688.332
$ws.Range("J65").Value = 57688.332
$ws.Range("L65").Value = 173064.996
$ws.Range("N65").Value = -179928.996

$ws.Range("H70").Value = 6867.1816
$ws.Range("I70").Value = 5814.636
$ws.Range("K70").Value = 5814.636
$ws.Range("M70").Value = -5544.636

$ws.Range("H73").Value = 6867.1816
$ws.Range("I73").Value = 5814.636
$ws.Range("K73").Value = 5814.636
$ws.Range("M73").Value = -4878.636

$ws.Range("H74").Value = 33420.332
$ws.Range("J74").Value = 33420.332
$ws.Range("L74").Value = 33420.332
$ws.Range("N74").Value = -35292.332

$ws.Range("H77").Value = 33420.332
$ws.Range("J77").Value = 33420.332
$ws.Range("L77").Value = 100260.996
$ws.Range("N77").Value = -109620.996

$ws.Range("H80").Value = 1399.8
$ws.Range("I80").Value = 1335
$ws.Range("J80").Value = 1497
$ws.Range("K80").Value = 1335
$ws.Range("L80").Value = 1497
$ws.Range("M80").Value = -337
$ws.Range("N80").Value = -3493

$ws.Range("H83").Value = 1399.8
$ws.Range("I83").Value = 1335
$ws.Range("J83").Value = 1497
$ws.Range("K83").Value = 6675
$ws.Range("L83").Value = 7485
$ws.Range("M83").Value = -1683
$ws.Range("N83").Value = -17469

$ws.Range("H141").Value = 89998
$ws.Range("J141").Value = 89998
$ws.Range("L141").Value = 89998
$ws.Range("N141").Value = -100358

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2376.75
$ws.Range("I7").Value = 2376.75
$ws.Range("K7").Value = 2376.75
$ws.Range("M7").Value = -2264.75

$ws.Range("H22").Value = 16700
$ws.Range("I22").Value = 28350
$ws.Range("K22").Value = 28350
$ws.Range("M22").Value = -28055

$ws.Range("H27").Value = 16700
$ws.Range("I27").Value = 28350
$ws.Range("K27").Value = 28350
$ws.Range("M27").Value = -28243

$ws.Range("H55").Value = 1157.9445
$ws.Range("I55").Value = 1094
$ws.Range("K55").Value = 1094
$ws.Range("M55").Value = -921

$ws.Range("H126").Value = 2376.75
$ws.Range("I126").Value = 2376.75
$ws.Range("K126").Value = 7130.25
$ws.Range("M126").Value = -4660.25

$ws.Range("H141").Value = 129985
$ws.Range("J141").Value = 129985
$ws.Range("L141").Value = 129985
$ws.Range("N141").Value = -140345

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6917.7646
$ws.Range("I136").Value = 7365.5347
$ws.Range("K136").Value = 22096.6041
$ws.Range("M136").Value = -19546.6041
